# Updated symbol list on Tue Jan  3 15:47:17 UTC 2023 with GitHub Actions
# Refresh Price (column D) and Volume(1h) (column E) figures for the
# crypto rows that moved since the last snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = "245.18"
    "E2" = "-0.67%"
    "D3" = "29.04"
    "E3" = "-1.34%"
    "D4" = "5.252"
    "E4" = "1.08%"
    "D5" = "0.05713"
    "E5" = "0.03%"
    "D6" = "6.609"
    "E6" = "0.44%"
    "D7" = "3.174"
    "E7" = "3.71%"
    "D8" = "0.8553"
    "E8" = "-0.38%"
    "D9" = "0.8565"
    "E9" = "-2.79%"
    "E10" = "0.06%"
    "D11" = "0.07040"
    "E11" = "-0.62%"
    "D12" = "0.03170"
    "E12" = "10.26%"
    "D13" = "0.09288"
    "E13" = "-1.02%"
    "D14" = "0.001525"
    "E14" = "0.26%"
    "D15" = "0.0005951"
    "E15" = "-0.54%"
    "D16" = "0.006079"
    "E16" = "-1.02%"
    "D17" = "3.509"
    "E17" = "0.82%"
    "E18" = "-4.31%"
    "D20" = "0.03329"
    "E20" = "0.73%"
    "D21" = "0.1277"
    "E21" = "-1.74%"
    "D22" = "3.484"
    "E22" = "0.49%"
    "D23" = "0.04136"
    "E23" = "-0.46%"
    "D24" = "0.1329"
    "E24" = "-3.58%"
    "E25" = "0.08%"
    "D26" = "0.004147"
    "E26" = "-17.97%"
    "D27" = "0.0001200"
    "E27" = "-0.72%"
    "D28" = "0.0001449"
    "E28" = "-25.24%"
    "D40" = "0.03835"
    "E40" = "2.17%"
    "D41" = "0.1066"
    "E41" = "-0.68%"
    "D42" = "0.002416"
    "E42" = "15.19%"
    "D43" = "0.002949"
    "E43" = "-47.83%"
    "D44" = "0.009407"
    "E44" = "-5.65%"
    "D45" = "0.00005277"
    "E45" = "2.91%"
    "D46" = "0.00000000750"
    "E46" = "0.13%"
    "D47" = "0.08991"
    "E47" = "26.75%"
    "D48" = "0.002452"
    "E48" = "-5.22%"
    "D49" = "0.00002100"
    "E49" = "0.13%"
    "D50" = "0.0002000"
    "E50" = "0.13%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force text storage so values such as "0.0002000" keep their
    # trailing zeros instead of being auto-parsed into numbers/percentages.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    # Restore the default (unstyled) cell style so the only change
    # recorded is the text content, matching the source diff.
    $cell.Style = "Normal"
}
